# Tasks have absolute dates
#
# Inserts two new columns ("Start Date", "Finish Date") right before the
# existing "Predecessors" column on every sheet, shifting the old F/G
# (Predecessors/Successors) to H/I, and fills the new columns with the
# absolute calendar dates corresponding to each task's relative Start
# Time / Finish Time (columns D/E), anchored at a 2023-01-08 project
# start date.

$wb = $excel.ActiveWorkbook

# Start Date (col F) / Finish Date (col G) literal values per sheet,
# one pair per data row (rows 2-6), pre-computed from base date
# 2023-01-08 + the Start Time / Finish Time (days) already in D/E.
$sheetDates = @{
    "P1_Constrained"    = @(
        @("08-01-2023", "08-01-2023"),
        @("08-01-2023", "29-03-2023"),
        @("29-03-2023", "17-06-2023"),
        @("17-06-2023", "05-09-2023"),
        @("05-09-2023", "05-09-2023")
    )
    "P1_notConstrained" = @(
        @("08-01-2023", "08-01-2023"),
        @("08-01-2023", "29-03-2023"),
        @("08-01-2023", "29-03-2023"),
        @("08-01-2023", "29-03-2023"),
        @("29-03-2023", "29-03-2023")
    )
    "P2_Constrained"    = @(
        @("18-01-2023", "18-01-2023"),
        @("18-01-2023", "08-04-2023"),
        @("08-04-2023", "27-06-2023"),
        @("27-06-2023", "15-09-2023"),
        @("15-09-2023", "15-09-2023")
    )
    "P2_notConstrained" = @(
        @("18-01-2023", "18-01-2023"),
        @("18-01-2023", "08-04-2023"),
        @("18-01-2023", "08-04-2023"),
        @("18-01-2023", "08-04-2023"),
        @("08-04-2023", "08-04-2023")
    )
    "P3_Constrained"    = @(
        @("28-01-2023", "28-01-2023"),
        @("28-01-2023", "18-04-2023"),
        @("18-04-2023", "07-07-2023"),
        @("07-07-2023", "25-09-2023"),
        @("25-09-2023", "25-09-2023")
    )
    "P3_notConstrained" = @(
        @("28-01-2023", "28-01-2023"),
        @("28-01-2023", "18-04-2023"),
        @("28-01-2023", "18-04-2023"),
        @("28-01-2023", "18-04-2023"),
        @("18-04-2023", "18-04-2023")
    )
}

foreach ($ws in $wb.Worksheets) {
    # Shift the existing Predecessors/Successors columns (F/G) two slots
    # to the right (-> H/I), opening up F/G for the new date columns.
    $ws.Range("F:G").Insert()

    $ws.Range("F1").Value = "Start Date"
    $ws.Range("G1").Value = "Finish Date"

    $dates = $sheetDates[$ws.Name]
    for ($i = 0; $i -lt $dates.Length; $i++) {
        $row = $i + 2
        $pair = $dates[$i]
        $ws.Range("F" + $row).Value = "'" + $pair[0]
        $ws.Range("G" + $row).Value = "'" + $pair[1]
    }
}
